$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Direct text assignments for Coin/Link/Volume columns (never parse as numbers).
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('E6').Value = '  +3.02%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -0.35%  '
$ws.Range('E9').Value = '  +5.49%  '
$ws.Range('E10').Value = '  -0.96%  '
$ws.Range('E11').Value = '  +0.35%  '
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('E13').Value = '  -0.01%  '
$ws.Range('E14').Value = '  +1.30%  '
$ws.Range('E15').Value = '  +0.33%  '
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('E17').Value = '  +0.85%  '
$ws.Range('E18').Value = '  -0.46%  '
$ws.Range('E19').Value = '  -0.90%  '
$ws.Range('E20').Value = '  +2.97%  '
$ws.Range('E21').Value = '  -1.90%  '
$ws.Range('E22').Value = '  -1.41%  '
$ws.Range('E23').Value = '  -3.02%  '
$ws.Range('E24').Value = '  +2.19%  '
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('E26').Value = '  -3.07%  '
$ws.Range('E27').Value = '  +4.26%  '
$ws.Range('E30').Value = '  -3.15%  '
$ws.Range('E31').Value = '  +0.39%  '
$ws.Range('E32').Value = '  -1.67%  '
$ws.Range('E33').Value = '  +0.25%  '
$ws.Range('E34').Value = '  -0.95%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('E36').Value = '  +1.41%  '
$ws.Range('E37').Value = '  +2.77%  '
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('E39').Value = '  +0.58%  '
$ws.Range('E40').Value = '  -2.40%  '
$ws.Range('E41').Value = '  -0.84%  '
$ws.Range('E42').Value = '  +4.62%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('E43').Value = '  -0.53%  '
$ws.Range('B44').Value = 'WhiteBITCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('E47').Value = '  -0.13%  '
$ws.Range('E48').Value = '  +0.88%  '
$ws.Range('E49').Value = '  -1.47%  '
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('E51').Value = '  +0.00%  '

# Price column (D) values that look like plain numbers must stay text, matching the
# original inline-string cells. Stage each one through a Text-formatted scratch cell and
# Copy/PasteSpecial(xlPasteValues) so the destination cell keeps its default (General)
# style while still receiving the literal text -- a direct .Value assignment would let
# Excel auto-convert numeric-looking strings ("75.00" -> 75, "0.0000191" -> 1.91E-05, etc).
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = '68.396.00'
$scratch.Copy() | Out-Null
$ws.Range('D2').PasteSpecial(-4163)
$scratch.Value = '2.650.29'
$scratch.Copy() | Out-Null
$ws.Range('D3').PasteSpecial(-4163)
$scratch.Value = '598.19'
$scratch.Copy() | Out-Null
$ws.Range('D5').PasteSpecial(-4163)
$scratch.Value = '159.71'
$scratch.Copy() | Out-Null
$ws.Range('D6').PasteSpecial(-4163)
$scratch.Value = '0.352'
$scratch.Copy() | Out-Null
$ws.Range('D12').PasteSpecial(-4163)
$scratch.Value = '28.11'
$scratch.Copy() | Out-Null
$ws.Range('D13').PasteSpecial(-4163)
$scratch.Value = '0.0000191'
$scratch.Copy() | Out-Null
$ws.Range('D14').PasteSpecial(-4163)
$scratch.Value = '3.131.28'
$scratch.Copy() | Out-Null
$ws.Range('D15').PasteSpecial(-4163)
$scratch.Value = '68.262.50'
$scratch.Copy() | Out-Null
$ws.Range('D16').PasteSpecial(-4163)
$scratch.Value = '2.663.99'
$scratch.Copy() | Out-Null
$ws.Range('D17').PasteSpecial(-4163)
$scratch.Value = '11.43'
$scratch.Copy() | Out-Null
$ws.Range('D18').PasteSpecial(-4163)
$scratch.Value = '363.74'
$scratch.Copy() | Out-Null
$ws.Range('D19').PasteSpecial(-4163)
$scratch.Value = '7.32'
$scratch.Copy() | Out-Null
$ws.Range('D21').PasteSpecial(-4163)
$scratch.Value = '4.82'
$scratch.Copy() | Out-Null
$ws.Range('D22').PasteSpecial(-4163)
$scratch.Value = '2.07'
$scratch.Copy() | Out-Null
$ws.Range('D23').PasteSpecial(-4163)
$scratch.Value = '75.00'
$scratch.Copy() | Out-Null
$ws.Range('D24').PasteSpecial(-4163)
$scratch.Value = '558.52'
$scratch.Copy() | Out-Null
$ws.Range('D30').PasteSpecial(-4163)
$scratch.Value = '159.78'
$scratch.Copy() | Out-Null
$ws.Range('D38').PasteSpecial(-4163)
$scratch.Value = '2.65'
$scratch.Copy() | Out-Null
$ws.Range('D43').PasteSpecial(-4163)
$scratch.Value = '17.80'
$scratch.Copy() | Out-Null
$ws.Range('D44').PasteSpecial(-4163)
$scratch.Value = '158.11'
$scratch.Copy() | Out-Null
$ws.Range('D46').PasteSpecial(-4163)
$scratch.Value = '22.27'
$scratch.Copy() | Out-Null
$ws.Range('D48').PasteSpecial(-4163)
$scratch.Value = '0.617'
$scratch.Copy() | Out-Null
$ws.Range('D51').PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = 0
